$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New rows being appended are exact duplicates of existing rows 3 and 2
# (row 4 <- row 3 content, row 5 <- row 2 content), matching the diff.
# Note: "Krishnappa Gowtham" is followed by a non-breaking space (U+00A0)
# in the original workbook, not a regular space - reuse that exact
# character so the duplicated cells are byte-identical to rows 3 / 2.
$nbsp = [char]0x00A0
$batsman = "Krishnappa Gowtham" + $nbsp
$row4 = @(" Abu Dhabi", " October 01 2020", "Mumbai won by 48 runs", "Kings XI Punjab", "Mumbai Indians", $batsman, "22", "13", "2", "1", "169.23")
$row5 = @(" Dubai (DSC)", " September 20 2020", "Match tied (Capitals won the one-over eliminator)", "Kings XI Punjab", "Delhi Capitals", $batsman, "20", "14", "1", "1", "142.85")

# Force the new cells to be stored as text (matching t="str" in the source
# sheet) instead of letting Excel auto-detect numeric-looking strings
# ("22", "13", "169.23", ...) as numbers.
$newRange = $ws.Range("A4:K5")
$newRange.NumberFormat = "@"

for ($c = 1; $c -le 11; $c++) {
    $ws.Cells.Item(4, $c).Value = $row4[$c - 1]
    $ws.Cells.Item(5, $c).Value = $row5[$c - 1]
}

# Restore default styling on the new cells so no lingering custom style
# index is left attached to them.
$newRange.ClearFormats()
